# The deck ships with two theme parts:
#   ppt/theme/theme1.xml -> bound to the (only) slide master, clrScheme "Integral"
#   ppt/theme/theme2.xml -> bound to the notes master,       clrScheme "Office"
# The authored change swaps the two themes' contents: the slide master
# (and therefore every slide) switches from the green/gold "Integral"
# palette over to the stock Office palette, while the notes master picks
# up the old "Integral" palette.
#
# The slide-facing theme is reachable and editable through
# Slide.ThemeColorScheme (a live view over the slide master's <a:clrScheme>),
# so drive the 12 theme colour slots there with the target ("Office Theme")
# RGB values. Order is fixed by PowerPoint's theme colour index:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2,
#   5 accent1, 6 accent2, 7 accent3, 8 accent4, 9 accent5, 10 accent6,
#   11 hlink, 12 folHlink
# Touching slide 1 is enough -- the colour scheme lives on the shared
# slide master/theme part, so every slide picks the new values up.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Item(1).RGB  = 0x000000   # dk1
$tcs.Item(2).RGB  = 0xFFFFFF   # lt1
$tcs.Item(3).RGB  = 0x6A5444   # dk2      -> 44546A
$tcs.Item(4).RGB  = 0xE6E6E7   # lt2      -> E7E6E6
$tcs.Item(5).RGB  = 0xD59B5B   # accent1  -> 5B9BD5
$tcs.Item(6).RGB  = 0x317DED   # accent2  -> ED7D31
$tcs.Item(7).RGB  = 0xA5A5A5   # accent3  -> A5A5A5
$tcs.Item(8).RGB  = 0x00C0FF   # accent4  -> FFC000
$tcs.Item(9).RGB  = 0xC47244   # accent5  -> 4472C4
$tcs.Item(10).RGB = 0x47AD70   # accent6  -> 70AD47
$tcs.Item(11).RGB = 0xC16305   # hlink    -> 0563C1
$tcs.Item(12).RGB = 0x724F95   # folHlink -> 954F72
